# UtilitarioCadastroTestCase addition: the "pkPessoa" header column (B1)
# is untouched content-wise; the real data edit is cell I2, whose label
# grows from "@TESTE" to "@TESTE - @teste" so the new utility test case
# can match on it. Leading "'" forces Excel to keep treating the text as
# an explicit string (preserves the existing quote-prefixed cell style
# instead of Excel minting a brand-new style without it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "'@TESTE - @teste"

# The sheet was left with the cursor resting on E7 when it was last saved.
$ws.Range("E7").Select() | Out-Null
